# 🔄 Actualización automática del tracker
# Applies the latest batch of results + newly scraped pending matches.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Fill in results ("resultado" / "profit") for matches that already
#    have a known outcome.
# ---------------------------------------------------------------------
$ws.Range("G18").Value = "Fallo"
$ws.Range("H18").Value = -1

$ws.Range("G20").Value = "Acierto"
$ws.Range("H20").Value = 0.73

$ws.Range("G32").Value = "Fallo"
$ws.Range("H32").Value = -1

# ---------------------------------------------------------------------
# 2) Append newly tracked matches (rows 43-53) that are still pending
#    (no resultado/profit yet).
# ---------------------------------------------------------------------
$newRows = @(
    @(43, 14644105, "2025-09-09", "Eliakim Coulibaly", "Leo Raquillet", "Gana Leo Raquillet", 3),
    @(44, 14644107, "2025-09-09", "Mikhail Kukushkin", "Enzo Couacaud", "Gana Enzo Couacaud", 2.1),
    @(45, 14644687, "2025-09-09", "Luca Castelnuovo", "Omar Jasika", "Gana Luca Castelnuovo", 2.38),
    @(46, 14644688, "2025-09-09", "Petr Bar Biryukov", "Tianhui Zhang", "Gana Tianhui Zhang", 2.63),
    @(47, 14644686, "2025-09-09", "Yuta Kikuchi", "Yu Hsiou Hsu", "Gana Yuta Kikuchi", 3.5),
    @(48, 14643824, "2025-09-09", "Imanol Lopez Morillo", "Daniel Rincon", "Gana Imanol Lopez Morillo", 3.25),
    @(49, 14643828, "2025-09-09", "Corentin Denolly", "Stefan Adrian Andreescu", "Gana Stefan Adrian Andreescu", 1.83),
    @(50, 14643829, "2025-09-09", "Mihai Alexandru Coman", "Cezar Cretu", "Gana Mihai Alexandru Coman", 8),
    @(51, 14643827, "2025-09-09", "Radu Mihai Papoe", "Gerard Campana Lee", "Gana Radu Mihai Papoe", 1.67),
    @(52, 14643826, "2025-09-09", "Emilien Demanet", "Nicolas Alvarez Varona", "Gana Emilien Demanet", 3.25),
    @(53, 14643825, "2025-09-09", "Stefan Palosi", "Mathys Erhard", "Gana Stefan Palosi", 3)
)

foreach ($row in $newRows) {
    $r = $row[0]

    $ws.Cells.Item($r, 1).Value = $row[1]

    # Keep the date as literal text (not auto-converted to a date serial).
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $row[2]

    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]

    # resultado / profit are still pending -> leave blank but make sure
    # the cells exist in the sheet (matches are not settled yet).
    $ws.Cells.Item($r, 7).NumberFormat = "General"
    $ws.Cells.Item($r, 8).NumberFormat = "General"
}
